$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "숙소에서 자가치료를 시도하는 이스티나. 이때 마침 메이가 숙소를 방문하게 된다.`n"
$ws.Range("D3").Value = "메이와 함께 온 건, 다름 아닌 우르수스 학생자치단의 또 다른 멤버, '나탈리야'였다.`n"

$ws.Rows.Item(2).AutoFit() | Out-Null
$ws.Rows.Item(3).AutoFit() | Out-Null
